$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (A2:F25) were reshuffled into a new row order (same row
# contents, different row positions). Write the full target grid back
# cell-by-cell since this runtime's Range.Value setter does not accept
# array objects.
$data = @(
    @(201, 9, 30, 15, 45, 30),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(101, 9, 30, 15, 60, 15),
    @(301, 6, 45, 30, 60, 45),
    @(401, 9, 48, 67, 75, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(501, 9, 52, 30, 75, 45),
    @(801, 3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(601, 9, 60, 67, 60, 42),
    @(502, 0, 4, 0, 0, 0),
    @(1, 0, 2, 2, 2, 2),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2),
    @(3, 0, 3, 3, 3, 3),
    @(1101, 0, 15, 30, 30, 0),
    @(402, 0, 0, 4, 0, 0),
    @(602, 0, 0, 4, 0, 9),
    @(702, 0, 0, 0, 4, 0),
    @(1002, 0, 0, 0, 0, 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}
